# "Evaluation update after bugs fixed"
# Updates the raw measurement data (columns C/D/E) in each of the four
# "Erkannt" tables on the sheet. The "Erkennungsrate [%]" column (F) holds
# formulas (=C/5, =D/5, =E/5 and the weighted-average row) and recalculates
# automatically from the new inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cosmetic: widen the sheet-tabs area slightly (tabRatio 305 -> 374 / 1000).
try { $excel.ActiveWindow.TabRatio = 0.374 } catch {}

# --- Table 1 (Accelerometer): rows 11-13 ---------------------------------
$ws.Range("C11").Value2 = 199
$ws.Range("D11").Value2 = 216
$ws.Range("E11").Value2 = 85

$ws.Range("C12").Value2 = 16
$ws.Range("D12").Value2 = 97
$ws.Range("E12").Value2 = 387

$ws.Range("C13").Value2 = 5
$ws.Range("D13").Value2 = 174
$ws.Range("E13").Value2 = 321

# --- Table 2 (Gyroskop): rows 19-21 --------------------------------------
$ws.Range("C19").Value2 = 151
$ws.Range("D19").Value2 = 215
$ws.Range("E19").Value2 = 134

$ws.Range("C20").Value2 = 2
$ws.Range("D20").Value2 = 80
$ws.Range("E20").Value2 = 418

$ws.Range("C21").Value2 = 0
$ws.Range("D21").Value2 = 129
$ws.Range("E21").Value2 = 371

# --- Table 3 (DCA (bias = 1.0)): rows 27-29 ------------------------------
$ws.Range("C27").Value2 = 187
$ws.Range("D27").Value2 = 245
$ws.Range("E27").Value2 = 68

$ws.Range("C28").Value2 = 195
$ws.Range("D28").Value2 = 146
$ws.Range("E28").Value2 = 159

$ws.Range("C29").Value2 = 11
$ws.Range("D29").Value2 = 43
$ws.Range("E29").Value2 = 446

# --- Table 4 (unnamed / 4th block): rows 35-37 ---------------------------
$ws.Range("C35").Value2 = 166
$ws.Range("D35").Value2 = 151
$ws.Range("E35").Value2 = 183

$ws.Range("C36").Value2 = 191
$ws.Range("D36").Value2 = 144
$ws.Range("E36").Value2 = 165

$ws.Range("C37").Value2 = 76
$ws.Range("D37").Value2 = 195
$ws.Range("E37").Value2 = 229

# Recalculate so the "Erkennungsrate" formulas pick up the new inputs.
$excel.Calculate()

# Restore the author's on-save cursor position/selection.
$null = $ws.Range("A19").Select()
$null = $ws.Range("D28").Select()
